# The deck currently uses the "Integral" design theme (ppt/theme/theme2.xml,
# linked from the SlideMaster) while a second, unused "Office Theme" theme
# part (ppt/theme/theme1.xml) only hangs off the Notes Master relationship.
#
# The authored edit swaps which theme carries which palette: the design that
# drives the slides switches from the Integral palette to the standard
# Office palette. We reproduce that by rewriting the live ThemeColorScheme
# (the 12-slot Dk1/Lt1/Dk2/Lt2/Accent1-6/Hyperlink/FollowedHyperlink scheme
# behind the presentation's one-and-only Design/SlideMaster) from the
# Integral colors to the Office Theme colors.
#
# Note PowerPoint's RGB long values are packed 0x00BBGGRR (blue high byte),
# i.e. the reverse byte order of the "RRGGBB" hex you see in the theme XML.

$p  = $ppt.ActivePresentation
$d  = $p.Designs.Item(1)
$sm = $d.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # dk1      -> #000000 (unchanged)
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1      -> #FFFFFF (unchanged)
$tcs.Item(3).RGB  = 0x6A5444   # dk2      -> #44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      -> #E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  -> #5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  -> #ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  -> #A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  -> #FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  -> #4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  -> #70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    -> #0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink -> #954F72
